# Math model B section edits (commit: "Start to working on Math model B text")
#
# Strategy: locate each target paragraph by its (unique) current text,
# then overwrite its content via Range.InsertXML with the exact OOXML
# the commit ends up with. This lets us control run-splitting / run-
# merging and w:pPr precisely, matching the canonical diff.

$d = $word.ActiveDocument

function Get-ParaByText($doc, $needle) {
    foreach ($p in $doc.Paragraphs) {
        if ($p.Range.Text -like $needle) {
            return $p
        }
    }
    return $null
}

# 1) "Seed" heading paragraph (pStyle 2): drop the stray <w:rPr> from <w:pPr>.
$p = Get-ParaByText $d "Seed*"
$p.Range.InsertXML('<w:p><w:pPr><w:pStyle w:val="2"/></w:pPr><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>Seed</w:t></w:r></w:p>')

# 2) "C 40%" -> "C 35%" split across three runs (same rPr each).
$p = Get-ParaByText $d "C 40%*"
$p.Range.InsertXML('<w:p><w:pPr><w:rPr><w:b/><w:bCs/><w:sz w:val="36"/><w:szCs w:val="40"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:hint="eastAsia"/><w:b/><w:bCs/><w:sz w:val="36"/><w:szCs w:val="40"/></w:rPr><w:t xml:space="preserve">C </w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/><w:b/><w:bCs/><w:sz w:val="36"/><w:szCs w:val="40"/></w:rPr><w:t>35</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/><w:b/><w:bCs/><w:sz w:val="36"/><w:szCs w:val="40"/></w:rPr><w:t>%</w:t></w:r></w:p>')

# 3) "D 15%" -> "D 20%" split across three runs.
$p = Get-ParaByText $d "D 15%*"
$p.Range.InsertXML('<w:p><w:pPr><w:rPr><w:b/><w:bCs/><w:sz w:val="36"/><w:szCs w:val="40"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:hint="eastAsia"/><w:b/><w:bCs/><w:sz w:val="36"/><w:szCs w:val="40"/></w:rPr><w:t xml:space="preserve">D </w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/><w:b/><w:bCs/><w:sz w:val="36"/><w:szCs w:val="40"/></w:rPr><w:t>20</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/><w:b/><w:bCs/><w:sz w:val="36"/><w:szCs w:val="40"/></w:rPr><w:t>%</w:t></w:r></w:p>')

# 4) "E 10%" -> "E 12%" split across three runs (new paragraph content).
$p = Get-ParaByText $d "E 10%*"
$p.Range.InsertXML('<w:p><w:pPr><w:rPr><w:b/><w:bCs/><w:sz w:val="36"/><w:szCs w:val="40"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:hint="eastAsia"/><w:b/><w:bCs/><w:sz w:val="36"/><w:szCs w:val="40"/></w:rPr><w:t>E 1</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/><w:b/><w:bCs/><w:sz w:val="36"/><w:szCs w:val="40"/></w:rPr><w:t>2</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/><w:b/><w:bCs/><w:sz w:val="36"/><w:szCs w:val="40"/></w:rPr><w:t>%</w:t></w:r></w:p>')

# 5) "F 5%" -> "F 3%" split across three runs (keeps lastRenderedPageBreak).
$p = Get-ParaByText $d "F 5%*"
$p.Range.InsertXML('<w:p><w:pPr><w:rPr><w:b/><w:bCs/><w:sz w:val="36"/><w:szCs w:val="40"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:hint="eastAsia"/><w:b/><w:bCs/><w:sz w:val="36"/><w:szCs w:val="40"/></w:rPr><w:lastRenderedPageBreak/><w:t xml:space="preserve">F </w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/><w:b/><w:bCs/><w:sz w:val="36"/><w:szCs w:val="40"/></w:rPr><w:t>3</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/><w:b/><w:bCs/><w:sz w:val="36"/><w:szCs w:val="40"/></w:rPr><w:t>%</w:t></w:r></w:p>')

# 6) "Begin with a random multi and a C on the board." - drop the <w:pPr> wrapper.
$p = Get-ParaByText $d "Begin with a random multi*"
$p.Range.InsertXML('<w:p><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>Begin with a random multi and a C on the board.</w:t></w:r></w:p>')

# 7) "Begins with a <mystery symbol> on the center of the board." - drop <w:pPr>.
$p = Get-ParaByText $d "Begins with a *"
$p.Range.InsertXML('<w:p><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t xml:space="preserve">Begins with a </w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>mystery symbol</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t xml:space="preserve"> on the center of the board.</w:t></w:r></w:p>')

# 8) "There are 2 types of Mystery symbol..." - drop <w:pPr>.
$p = Get-ParaByText $d "There are 2 types*"
$p.Range.InsertXML('<w:p><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t xml:space="preserve">There are 2 types of Mystery symbol. The weaker one being </w:t></w:r><w:r><w:t>&#8220;</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>?</w:t></w:r><w:r><w:t>&#8221;</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t xml:space="preserve">, and the stronger one being </w:t></w:r><w:r><w:t>&#8220;</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>???</w:t></w:r><w:r><w:t>&#8221;</w:t></w:r></w:p>')

# 9) "70%: A coin. (Min 50x)" - drop <w:pPr>.
$p = Get-ParaByText $d "70%: A coin*"
$p.Range.InsertXML('<w:p><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>70%: A coin. (Min 50x)</w:t></w:r></w:p>')

# 10) "8%: A collect. (If no pending Collect is on the board)" - drop <w:pPr>.
$p = Get-ParaByText $d "8%: A collect*"
$p.Range.InsertXML('<w:p><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>8%: A collect. (If no pending Collect is on the board)</w:t></w:r></w:p>')

# 11) merge the split "?" + "??" runs in the strong-symbol paragraph into "???".
$p = Get-ParaByText $d "*???*Symbol will contain one of the following:*"
$p.Range.InsertXML('<w:p><w:r><w:t>&#8220;</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>???</w:t></w:r><w:r><w:t>&#8221;</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t xml:space="preserve"> Symbol will contain one of the following:</w:t></w:r></w:p>')

# 12) "63%: A coin. (Min 250x)" - merge runs, drop <w:pPr>.
$p = Get-ParaByText $d "63*A coin*"
$p.Range.InsertXML('<w:p><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>63%: A coin. (Min 250x)</w:t></w:r></w:p>')

# 13) "26%: A multiplier (Min 5x)" - merge runs, drop <w:pPr>.
$p = Get-ParaByText $d "26*A multiplier*"
$p.Range.InsertXML('<w:p><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>26%: A multiplier (Min 5x)</w:t></w:r></w:p>')

# 14) "11%: A collect. (If no pending Collect is on the board)" - merge runs, drop <w:pPr>.
$p = Get-ParaByText $d "11*A collect*"
$p.Range.InsertXML('<w:p><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>11%: A collect. (If no pending Collect is on the board)</w:t></w:r></w:p>')

Write-Output "done"
